$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (logistic_embeddings)
$ws.Range("C5").Value = 0.642
$ws.Range("D5").Value = 0.738
$ws.Range("E5").Value = 0.763
$ws.Range("F5").Value = 0.803
$ws.Range("G5").Value = 0.705
$ws.Range("H5").Value = 0.71

# Row 7 (classical-best-embeddings -> classical-best-embed)
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.642
$ws.Range("E7").Value = 0.763
$ws.Range("F7").Value = 0.803
$ws.Range("G7").Value = 0.705
$ws.Range("H7").Value = 0.71

# Row 8 (BERT-base)
$ws.Range("C8").Value = 0.668
$ws.Range("D8").Value = 0.768
$ws.Range("E8").Value = 0.797
$ws.Range("F8").Value = 0.803
$ws.Range("G8").Value = 0.742
$ws.Range("H8").Value = 0.753

# Row 9 (BERT-base-nli)
$ws.Range("B9").Value = 0.543
$ws.Range("C9").Value = 0.63
$ws.Range("D9").Value = 0.747
$ws.Range("E9").Value = 0.772
$ws.Range("F9").Value = 0.776
$ws.Range("G9").Value = 0.6899999999999999
$ws.Range("H9").Value = 0.6899999999999999
